$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '62.640.19'
Set-TextValue "E2" '  -3.49%  '

Set-TextValue "D3" '3.051.52'
Set-TextValue "E3" '  -2.71%  '

Set-TextValue "E4" '  +0.06%  '

Set-TextValue "D5" '543.09'
Set-TextValue "E5" '  -3.83%  '

Set-TextValue "D6" '133.88'
Set-TextValue "E6" '  -9.22%  '

Set-TextValue "D7" '1.00'
Set-TextValue "E7" '  +0.09%  '

Set-TextValue "D8" '3.042.43'
Set-TextValue "E8" '  -2.74%  '

Set-TextValue "D9" '0.487'
Set-TextValue "E9" '  -1.89%  '

Set-TextValue "D10" '6.46'
Set-TextValue "E10" '  -6.63%  '

Set-TextValue "E11" '  -1.12%  '

Set-TextValue "E12" '  -0.96%  '

Set-TextValue "D13" '34.80'
Set-TextValue "E13" '  -2.71%  '

Set-TextValue "D14" '0.0000214'
Set-TextValue "E14" '  -3.19%  '

Set-TextValue "D15" '3.547.08'
Set-TextValue "E15" '  -2.36%  '

Set-TextValue "D16" '62.633.38'
Set-TextValue "E16" '  -3.47%  '

Set-TextValue "E17" '  -1.34%  '

Set-TextValue "D18" '3.050.37'
Set-TextValue "E18" '  -2.71%  '

Set-TextValue "D19" '6.58'
Set-TextValue "E19" '  -1.91%  '

Set-TextValue "D20" '480.87'
Set-TextValue "E20" '  -8.38%  '

Set-TextValue "D21" '13.34'
Set-TextValue "E21" '  -3.23%  '

Set-TextValue "D22" '0.695'
Set-TextValue "E22" '  -0.62%  '

Set-TextValue "D23" '7.01'
Set-TextValue "E23" '  -4.95%  '

Set-TextValue "D24" '77.24'
Set-TextValue "E24" '  -1.59%  '

Set-TextValue "D25" '12.14'
Set-TextValue "E25" '  -4.37%  '

Set-TextValue "E26" '  +0.18%  '

Set-TextValue "D27" '2.70'
Set-TextValue "E27" '  -3.29%  '

Set-TextValue "D28" '8.21'
Set-TextValue "E28" '  -4.63%  '

Set-TextValue "D29" '0.997'
Set-TextValue "E29" '  -0.23%  '

Set-TextValue "D30" '1.92'
Set-TextValue "E30" '  -9.13%  '

Set-TextValue "D31" '26.10'
Set-TextValue "E31" '  -0.04%  '

Set-TextValue "D32" '1.12'
Set-TextValue "E32" '  -3.06%  '

Set-TextValue "B33" 'Stacks'
Set-TextValue "C33" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D33" '2.48'
Set-TextValue "E33" '  -6.43%  '

Set-TextValue "B34" 'OKB'
Set-TextValue "C34" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D34" '59.80'
Set-TextValue "E34" '  +13.28%  '

Set-TextValue "D35" '507.52'
Set-TextValue "E35" '  -8.62%  '

Set-TextValue "D36" '5.90'
Set-TextValue "E36" '  -1.71%  '

Set-TextValue "D37" '5.10'
Set-TextValue "E37" '  -5.24%  '

Set-TextValue "D38" '0.0397'
Set-TextValue "E38" '  -9.29%  '

Set-TextValue "D39" '3.063.69'
Set-TextValue "E39" '  -0.13%  '

Set-TextValue "D40" '0.0784'
Set-TextValue "E40" '  -3.46%  '

Set-TextValue "D41" '0.116'
Set-TextValue "E41" '  -3.56%  '

Set-TextValue "D42" '8.01'
Set-TextValue "E42" '  -2.52%  '

Set-TextValue "D43" '2.60'
Set-TextValue "E43" '  -8.42%  '

Set-TextValue "D44" '0.251'
Set-TextValue "E44" '  -1.57%  '

Set-TextValue "D46" '2.03'
Set-TextValue "E46" '  -5.98%  '

Set-TextValue "D47" '120.53'
Set-TextValue "E47" '  +1.25%  '

Set-TextValue "D48" '24.18'
Set-TextValue "E48" '  -2.92%  '

Set-TextValue "D49" '0.106'
Set-TextValue "E49" '  -1.68%  '

Set-TextValue "D50" '0.0₃0494'
Set-TextValue "E50" '  -5.33%  '

Set-TextValue "D51" '2.34'
Set-TextValue "E51" '  +61.91%  '
